$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1.4
$ws.Range("E3").Value = 1.29
$ws.Range("F3").Value = 1.2
$ws.Range("E4").Value = 1.22
$ws.Range("C5").Value = 1.39
$ws.Range("G5").Value = 0.76
$ws.Range("C6").Value = 1.49
$ws.Range("E6").Value = 1.32
